$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value of B2: "Pro-Lot Track (Lot Track)" -> "Pro-SYDATA1 (Lot track)"
$ws.Range("B2").Value = "Pro-SYDATA1 (Lot track)"

# Remove the bold Arial header formatting from A1:F1 (revert to default/Normal style)
$ws.Range("A1:F1").Style = "Normal"

# Update the sheet selection to cover the whole data range (A1:K2)
$ws.Activate()
$ws.Range("A1:K2").Select()
